# "Début d'un système de sons" - add a new reference link (row 23) for a
# StackOverflow answer about non-blocking playsound, matching the other
# hyperlinked URL rows already on the "Aide (Référence)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://stackoverflow.com/questions/41421313/make-playsound-non-blocking"

# New row right after the last existing reference (row 22).
$cell = $ws.Range("A23")
$cell.Value = $newUrl

# Turn the cell into a hyperlink, like every other reference row.
$ws.Hyperlinks.Add($cell, $newUrl) | Out-Null

# Hyperlinks.Add() tends to stamp its own ad-hoc formatting; re-apply the
# same built-in hyperlink style used by the other reference cells so A23
# matches its neighbours (s="2").
$cell.Style = "Lien hypertexte"

# Leave the selection where the author ended up after entering the new
# link (a couple of rows below the freshly typed entry).
$ws.Range("A25").Select() | Out-Null
